$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells with same style as the existing header row (copy format from AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96  # AD
    $ws.Cells.Item($r, 31).Value = 66  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
